# Hoàn thiện Ngoại Trú
# Finalize the outpatient ("Ngoại Trú") receiving test-case row on the
# "Data" sheet (row 2) with the updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Id changed
$ws.Range("A2").Value = 3018

# IdCardNo changed
$ws.Range("E2").Value = 46200608018

# InsBenefitType / InsBenefitRatio changed
$ws.Range("AM2").Value = 1
$ws.Range("AN2").Value = 0

# InsCheckedMessage ("Thẻ BHYT hợp lệ") and MedServiceId (4803) are no
# longer part of the test data - clear them out entirely.
$ws.Range("AQ2").ClearContents()
$ws.Range("AS2").ClearContents()

# WardUnitId is now populated.
$ws.Range("AT2").Value = 149

# Leave the cursor where the edit last happened.
$null = $ws.Range("AQ2").Select()
